# Edit script: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# - Adds a new worker-period row (2509) to the EC table
# - Updates VALOR MORA and Cant. Periodos totals accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row before row 18 (shifts the signature block rows down by one)
$ws.Rows("18:18").Insert()

# 2. Populate the new row 18 with the new period's data (mirrors rows 16/17)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "7921698"
$ws.Range("D18").Value = "LUIS CARLOS PUELLO GONZALEZ"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# 3. Match formatting of the row above (font, number format, borders, fill)
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4123)
$excel.CutCopyMode = 0

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "7921698"
$ws.Range("D18").Value = "LUIS CARLOS PUELLO GONZALEZ"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Apply full box borders (thin) to the new row, matching the rest of the table
$rowRange = $ws.Range("B18:J18")
$rowRange.Borders.Item(7).LineStyle = 1
$rowRange.Borders.Item(7).Weight = 2
$rowRange.Borders.Item(8).LineStyle = 1
$rowRange.Borders.Item(8).Weight = 2
$rowRange.Borders.Item(9).LineStyle = 1
$rowRange.Borders.Item(9).Weight = 2
$rowRange.Borders.Item(10).LineStyle = 1
$rowRange.Borders.Item(10).Weight = 2
$rowRange.Borders.Item(11).LineStyle = 1
$rowRange.Borders.Item(11).Weight = 2
$rowRange.Borders.Item(12).LineStyle = 1
$rowRange.Borders.Item(12).Weight = 2

# 4. Update the summary values: Valor Mora and Cant. Periodos
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3

Write-Host "Edit complete"
